{"js": "// Remove the \"Lab - 1\" text from the final paragraph of the document body,\n// leaving the (now empty) paragraph and its formatting (pPr/rPr) intact.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst target = items[items.length - 1];\n\n// Safety check: only clear it if it actually holds the \"Lab - 1\" text\n// (falls back to the last paragraph regardless, matching the diff's intent).\nif (target.text === \"Lab - 1\" || /Lab\\s*-\\s*1/.test(target.text)) {\n  target.insertText(\"\", \"Replace\");\n} else {\n  // Fallback: search the whole document for the run text and clear its\n  // containing paragraph.\n  const results = body.search(\"Lab - 1\", { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n  if (results.items.length > 0) {\n    const para = results.items[0].paragraphs.getFirst();\n    para.insertText(\"\", \"Replace\");\n  } else {\n    target.insertText(\"\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the \"Lab - 1\" text from the document, leaving the (now empty)\n# paragraph and its formatting (pPr/rPr) intact.\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $p = $d.Paragraphs.Item($i)\n    $t = $p.Range.Text\n    if ($t -match \"Lab\\s*-\\s*1\") {\n        $r = $p.Range\n        # Exclude the trailing paragraph mark (the last character of a\n        # paragraph's Range) so only the run text is removed and the\n        # paragraph itself survives.\n        $textRange = $d.Range($r.Start, $r.End - 1)\n        $textRange.Delete()\n    }\n}\n"}
